# dev-template.xlsx: wire up the "Prix Total" column (H) in the quote-line
# table so it auto-computes Qte (F) * PU Net (G) for every line item row,
# and leave the view where the user was last working (bottom of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 20-39 are the blank line-item rows of the quote table
# (A=Reference, B:E=Designation, F=Qte, G=PU Net, H=Prix Total).
# Give H on each row a real formula instead of a blank styled cell.
for ($r = 20; $r -le 39; $r++) {
    $ws.Range("H$r").Formula = "=SUM(G$r*F$r)"
}

# Recalculate so the cached <v> values are correct.
$excel.CalculateFull()

# Scroll the sheet down so row 13 is at the top and move the active
# selection to M22, matching where the author left the view.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("M22").Select()
